$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 7.172092999999999
$ws.Range("H2").Value = 21.516279
$ws.Range("I2").Value = 0.07357387076805701
$ws.Range("J2").Value = 0.07357387076805699
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 2.08532
$ws.Range("N2").Value = 6.25596
$ws.Range("O2").Value = 0.01753772176136817
$ws.Range("P2").Value = 0.01753772176136816
$ws.Range("Q2").Value = 14.95610897476
$ws.Range("R2").Value = 134.60498077284
$ws.Range("S2").Value = 0.001290318074437043
$ws.Range("T2").Value = 0.001290318074437042

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 7.172092999999999
$ws.Range("H3").Value = 21.516279
$ws.Range("I3").Value = 0.07357387076805701
$ws.Range("J3").Value = 0.07357387076805699
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 101.898173
$ws.Range("N3").Value = 305.694519
$ws.Range("O3").Value = 0.8569724579756384
$ws.Range("P3").Value = 0.8569724579756383
$ws.Range("Q3").Value = 730.8231732860889
$ws.Range("R3").Value = 6577.4085595748
$ws.Range("S3").Value = 0.06305078087488378
$ws.Range("T3").Value = 0.06305078087488376

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 7.172092999999999
$ws.Range("H4").Value = 21.516279
$ws.Range("I4").Value = 0.07357387076805701
$ws.Range("J4").Value = 0.07357387076805699
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 14.921347
$ws.Range("N4").Value = 44.76404100000001
$ws.Range("O4").Value = 0.1254898202629935
$ws.Range("P4").Value = 0.1254898202629935
$ws.Range("Q4").Value = 107.017288369271
$ws.Range("R4").Value = 963.155595323439
$ws.Range("S4").Value = 0.009232771818736187
$ws.Range("T4").Value = 0.009232771818736182

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 89.72947699999999
$ws.Range("H5").Value = 269.188431
$ws.Range("I5").Value = 0.9204767624852804
$ws.Range("J5").Value = 0.9204767624852804
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 2.08532
$ws.Range("N5").Value = 6.25596
$ws.Range("O5").Value = 0.01753772176136817
$ws.Range("P5").Value = 0.01753772176136816
$ws.Range("Q5").Value = 187.11467297764
$ws.Range("R5").Value = 1684.03205679876
$ws.Range("S5").Value = 0.01614306534827182
$ws.Range("T5").Value = 0.01614306534827182

# Row 6
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 89.72947699999999
$ws.Range("H6").Value = 269.188431
$ws.Range("I6").Value = 0.9204767624852804
$ws.Range("J6").Value = 0.9204767624852804
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 101.898173
$ws.Range("N6").Value = 305.694519
$ws.Range("O6").Value = 0.8569724579756384
$ws.Range("P6").Value = 0.8569724579756383
$ws.Range("Q6").Value = 9143.269770545519
$ws.Range("R6").Value = 82289.42793490969
$ws.Range("S6").Value = 0.7888232336564687
$ws.Range("T6").Value = 0.7888232336564686

# Row 7
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 89.72947699999999
$ws.Range("H7").Value = 269.188431
$ws.Range("I7").Value = 0.9204767624852804
$ws.Range("J7").Value = 0.9204767624852804
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 14.921347
$ws.Range("N7").Value = 44.76404100000001
$ws.Range("O7").Value = 0.1254898202629935
$ws.Range("P7").Value = 0.1254898202629935
$ws.Range("Q7").Value = 1338.884662445519
$ws.Range("R7").Value = 12049.96196200967
$ws.Range("S7").Value = 0.11551046348054
$ws.Range("T7").Value = 0.11551046348054

# Row 8
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 0.5799533333333334
$ws.Range("H8").Value = 1.73986
$ws.Range("I8").Value = 0.005949366746662454
$ws.Range("J8").Value = 0.005949366746662453
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 2.08532
$ws.Range("N8").Value = 6.25596
$ws.Range("O8").Value = 0.01753772176136817
$ws.Range("P8").Value = 0.01753772176136816
$ws.Range("Q8").Value = 1.209388285066667
$ws.Range("R8").Value = 10.8844945656
$ws.Range("S8").Value = 0.0001043383386593023
$ws.Range("T8").Value = 0.0001043383386593022

# Row 9
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 0.5799533333333334
$ws.Range("H9").Value = 1.73986
$ws.Range("I9").Value = 0.005949366746662454
$ws.Range("J9").Value = 0.005949366746662453
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 101.898173
$ws.Range("N9").Value = 305.694519
$ws.Range("O9").Value = 0.8569724579756384
$ws.Range("P9").Value = 0.8569724579756383
$ws.Range("Q9").Value = 59.09618509192668
$ws.Range("R9").Value = 531.8656658273401
$ws.Range("S9").Value = 0.005098443444285851
$ws.Range("T9").Value = 0.005098443444285849

# Row 10
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 0.5799533333333334
$ws.Range("H10").Value = 1.73986
$ws.Range("I10").Value = 0.005949366746662454
$ws.Range("J10").Value = 0.005949366746662453
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 14.921347
$ws.Range("N10").Value = 44.76404100000001
$ws.Range("O10").Value = 0.1254898202629935
$ws.Range("P10").Value = 0.1254898202629935
$ws.Range("Q10").Value = 8.653684930473336
$ws.Range("R10").Value = 77.88316437426002
$ws.Range("S10").Value = 0.0007465849637173019
$ws.Range("T10").Value = 0.0007465849637173016
